# The barcode-scanning script used to only ever update the *last* row of
# the sheet ("read last one"). It has been fixed to read/append every
# scanned barcode, so the sheet now grows with one row per scanned code
# instead of just overwriting the final row.
#
# This applies the resulting data changes: row 3's quantity is corrected,
# row 4 now holds a different (previously-seen) barcode with an updated
# quantity, and four brand-new barcode rows (5-8) are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing row 3: quantity updated (2 -> 5)
$ws.Range("B3").Value = 5

# Existing row 4: barcode + quantity replaced
$ws.Range("A4").Value = 7610700949085
$ws.Range("B4").Value = 7

# Newly scanned barcodes appended below the previous data
$ws.Range("A5").Value = 3502110008091
$ws.Range("B5").Value = 20

$ws.Range("A6").Value = 5410013110002
$ws.Range("B6").Value = 5

$ws.Range("A7").Value = 5411028070480
$ws.Range("B7").Value = 90

$ws.Range("A8").Value = 5411188115472
$ws.Range("B8").Value = 63

# Reflect the active selection left behind after the batch of updates
$ws.Range("A2:B8").Select()
